# Docs: final presentation update
#
# Target edit (per the authoritative OOXML diff):
#   1. Slide 3 ("Context & Motivation" source-citation slide) loses its
#      "Slide Number Placeholder 1" shape (<p:sp> with ph type="sldNum")
#      entirely -- it is removed from the slide, not just emptied.
#   2. Two embedded Visio OLE graphicFrames (slide 7 and slide 12) have
#      their legacy VML fallback `spid` attribute incremented by one
#      (e.g. _x0000_s1038 -> _x0000_s1039). That attribute is purely an
#      internal legacy-VML bookkeeping id inside the <mc:Choice> fallback
#      of the OLE p:oleObj element; it is not surfaced anywhere on the
#      PowerPoint Shape/OLEFormat object model (confirmed against the
#      full Get-Member property list for the shape), so it cannot be
#      targeted from COM-interop script -- nothing else about those two
#      OLE objects changes (same r:id, same image, same size/position).

$p = $ppt.ActivePresentation

# --- Slide 3: remove the Slide Number Placeholder shape completely ----
# Placeholder shapes that are still defined on the layout get "reset" to
# an empty layout-inherited stub the first time Delete() is called (this
# mirrors how PowerPoint keeps "required" placeholders alive). Calling
# Delete() a second time on that freshly reset stub removes it for good,
# matching the target XML (the shape is gone, nothing left behind).
$slide3 = $p.Slides.Item(3)

$sldNumShape = $null
for ($i = 1; $i -le $slide3.Shapes.Count; $i++) {
    $candidate = $slide3.Shapes.Item($i)
    if ($candidate.Name -like "Slide Number Placeholder*") {
        $sldNumShape = $candidate
        break
    }
}

if ($sldNumShape -ne $null) {
    $sldNumShape.Delete()

    $sldNumShape2 = $null
    for ($i = 1; $i -le $slide3.Shapes.Count; $i++) {
        $candidate = $slide3.Shapes.Item($i)
        if ($candidate.Name -like "Slide Number Placeholder*") {
            $sldNumShape2 = $candidate
            break
        }
    }
    if ($sldNumShape2 -ne $null) {
        $sldNumShape2.Delete()
    }
}
